# Adds a new "區隔帳冊" (segregated ledger/book) column to the right of the
# existing data table on sheet "1091231", matching the formatting of the
# neighbouring columns, resizes it, and moves the AutoFilter (and the
# worksheet's hidden _FilterDatabase name that backs it) from A1:L1 to
# C1:M1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("1091231")

# --- 1. New header cell M1, formatted like the other header cells (C1:L1) ---
$ws.Range("L1").Copy()
$ws.Range("M1").PasteSpecial(-4122)              # xlPasteFormats
$ws.Range("M1").Value = "區隔帳冊"

# --- 2. New data cells M2:M3, formatted like the neighbouring K2:L3 data cells ---
$ws.Range("L2:L3").Copy()
$ws.Range("M2:M3").PasteSpecial(-4122)           # xlPasteFormats

$excel.CutCopyMode = $false

# --- 3. Give the new column an explicit width (matches the other sized columns) ---
$ws.Range("M1").ColumnWidth = 10.285714285714286

# --- 4. Re-point the AutoFilter so it spans C1:M1 instead of A1:L1 -----------
$ws.AutoFilterMode = $false
$ws.Range("C1:M1").AutoFilter()

# Keep the worksheet-scoped hidden _FilterDatabase name in sync with the
# AutoFilter range above (Excel keeps these two in lock-step).
$filterDbName = $wb.Names.Item("1091231!_FilterDatabase")
$filterDbName.RefersTo = "='1091231'!`$C`$1:`$M`$1"

# --- 5. Update the view so the new column is visible / selected -------------
$ws.Activate()
$ws.Range("I10").Select()
$excel.ActiveWindow.ScrollColumn = 7
